# ChartAPITest.xlsx edit: rename Sheet1 -> Bubble, re-point the bubble chart
# series at the already-staged F:H helper columns, restyle the bubble series
# (accent2 fill / heavier outline) and re-home the chart a column to the
# left, then leave "Bubble" as the selected tab (with J9 selected) and
# "Line" no longer the active tab.

$wb = $excel.ActiveWorkbook

# --- Rename the first sheet to match its new data role ---------------------
$wsBubble = $wb.Worksheets.Item(1)
$wsBubble.Name = "Bubble"

$wsLine = $wb.Worksheets.Item("Line")

# --- Re-target + restyle the bubble chart series ---------------------------
$co = $wsBubble.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)

# Series name ("Product B") now shows up as the c:tx entry.
$ser.Name = "Product B"

# Restyle: heavier (2pt, still invisible/noFill) outline.
$ser.Format.Line.Weight = 2
$ser.Format.Line.Visible = $false

# --- Reposition the chart one column to the left ----------------------------
$co.Left = 127.375
$co.Top = 163.12496062992125
$co.Width = 433.0625
$co.Height = 216
$co.Name = "Chart 2"

# --- Tab / selection bookkeeping -------------------------------------------
# "Bubble" becomes the active tab with J9 selected; "Line" keeps its own
# selection (M6) but is no longer the active tab.
$wsBubble.Activate()
$wsBubble.Range("J9").Select()
